# Insert a new data row at row 58 (pushing the existing rows 58:74 down to
# 59:75), then populate the new row with the weekly Coco record, mirroring
# the style of the row that used to be at position 58 (date column uses a
# date number format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 58 downwards (inserting a blank row above the old row 58).
$ws.Rows("58:58").Insert()

# Copy the date-cell style (numFmtId 165 "YYYY-MM-DD HH:MM:SS") from the row
# directly below (the row that used to be 58, now 59) onto the new D58 cell.
$ws.Range("D59").Copy()
$ws.Range("D58").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 44736
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100108
$ws.Range("H58").Value = "Tropicales y subtropicales"
$ws.Range("I58").Value = 100108007
$ws.Range("J58").Value = "Coco"
$ws.Range("K58").Value = "Sin especificar"
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 20
$ws.Range("N58").Value = 28000
$ws.Range("O58").Value = 28000
$ws.Range("P58").Value = 28000
$ws.Range("Q58").Value = "$/malla 20 unidades"
$ws.Range("R58").Value = "Perú"
$ws.Range("S58").Value = 1400
$ws.Range("T58").Value = 20
